# Update "Datos actualizados" timestamp cell
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

$ws.Range("A1").Value = "Datos actualizados a 18 de Mayo de 2020 a las 19:05"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1534453
$ws.Range("C4").Value = 6789
$ws.Range("D4").Value = 347253
$ws.Range("E4").Value = 1095932
$ws.Range("G4").Value = 290
$ws.Range("H4").Value = 91268

# Row 11 - Alemania
$ws.Range("B11").Value = 177182
$ws.Range("C11").Value = 531
$ws.Range("E11").Value = 14497
$ws.Range("G11").Value = 36
$ws.Range("H11").Value = 8085

# Row 14 - India
$ws.Range("B14").Value = 100293
$ws.Range("C14").Value = 4595
$ws.Range("E14").Value = 58229
$ws.Range("G14").Value = 130
$ws.Range("H14").Value = 3155

# Row 135 - Nepal
$ws.Range("B135").Value = 375
$ws.Range("C135").Value = 80
$ws.Range("E135").Value = 337

# Row 154 - Birmania
$ws.Range("B154").Value = 188
$ws.Range("C154").Value = 4
$ws.Range("D154").Value = 101
$ws.Range("E154").Value = 81
